# Update countries & provincias Spain
# Applies the 25-Jun-2020 refresh (21:53 -> 23:10) to the "Pais" sheet:
#  - updates the "Datos actualizados..." timestamp
#  - updates Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes for the countries whose figures
#    moved
#  - re-sorts a few countries whose rank changed (Egipto overtook
#    Bielorrusia & Belgica; Fiyi/Dominica, Groenlandia/Islas Malvinas and
#    Seychelles/Montserrat swapped rank order further down the table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $pais, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($r, 1).Value = $pais
    $ws.Cells.Item($r, 2).Value = $casosTotales
    $ws.Cells.Item($r, 3).Value = $nuevosCasos
    $ws.Cells.Item($r, 4).Value = $casosActivos
    $ws.Cells.Item($r, 5).Value = $recuperados
    $ws.Cells.Item($r, 6).Value = $casosCriticos
    $ws.Cells.Item($r, 7).Value = $muertesHoy
    $ws.Cells.Item($r, 8).Value = $muertes
}

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 23:10"

# Figure refreshes (country keeps its row, only the metrics change)
Set-Row 4   "Estados Unidos" 2493141 30587 1045103 1323249 0 508 124789
Set-Row 7   "India"          491170  18185 285671  190191  0 401 15308
Set-Row 10  "Peru"           268602  3913  156074  103767  0 175 8761
Set-Row 15  "Alemania"       193663  409   176800  7851    0 9   9012
Set-Row 50  "Barein"         24081   511   18501   5509    0 2   71
Set-Row 108 "Mali"           2007    2     1354    540     0 1   113
Set-Row 140 "Ruanda"         850     20    385     463     0 0   2
Set-Row 156 "Montenegro"     414     25    315     90      0 0   9

# Egipto overtakes Bielorrusia and Belgica (rows 27-29 shift down the
# ranking, each now carries the figures of the country that used to sit
# one place above it, and Egipto gets its refreshed totals)
Set-Row 27  "Egipto"      61130 1569 16338 42259 0 83 2533
Set-Row 28  "Belgica"     61007 109  16890 34391 0 4  9726
Set-Row 29  "Bielorrusia" 60382 437  41448 18567 0 5  367

# Rank swaps further down the table (figures unchanged, only order)
Set-Row 202 "Fiyi"      18 0 18 0 0 0 0
Set-Row 203 "Dominica"  18 0 18 0 0 0 0

Set-Row 208 "Groenlandia"   13 0 13 0 0 0 0
Set-Row 209 "Islas Malvinas" 13 0 13 0 0 0 0

Set-Row 211 "Seychelles" 11 0 11 0 0 0 0
Set-Row 212 "Montserrat" 11 0 10 0 0 0 1
